# "temp reader - first commit"
#
# Reworks the Sheet1 listing from a DIR/CAD/PDC "converted document" table
# into a table of numeric file ids, and drops the third data row entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the last data row (previously the "PDC" row, row 4) -----------
# This shifts rows 2:4 -> up by one and shrinks the used range to A1:D3.
$ws.Rows(4).Delete()

# --- Column A no longer carries a value on the data rows -------------------
# (the "type" header in A1 stays; A2:A3 become blank)
$ws.Range("A2:A3").ClearContents()

# --- Column B becomes a numeric file id instead of the literal text --------
$ws.Range("B2").Value = 100014616490
$ws.Range("B3").Value = 100033325020

# Integer display format for the header + the two id cells
$ws.Range("B1:B3").NumberFormat = "0"

# Widen column B to fit the longer id values (best-fit-style width)
$ws.Columns("B").ColumnWidth = 18.8

# --- Selection / view bookkeeping ------------------------------------------
$null = $ws.Range("C4").Select()

# --- Footer stamp added to the sheet ---------------------------------------
# PageSetup.CenterFooter writes into <oddFooter> with an implicit "&C" prefix
$ws.PageSetup.CenterFooter = "_x000D_&1#&`"Calibri`"&10&K000000 SLB-Private"
